$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2074
$ws.Range("K3").Value = 2010
$ws.Range("F4").Value = 1907
$ws.Range("J4").Value = 1805
$ws.Range("K4").Value = 423
$ws.Range("K5").Value = 136
$ws.Range("K6").Value = 2548
$ws.Range("F7").Value = 24100
$ws.Range("J7").Value = 29276
$ws.Range("K7").Value = 7191

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 138
$ws.Range("K3").Value = 138
$ws.Range("K7").Value = 484

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 81
$ws.Range("K7").Value = 288

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 38
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 61
$ws.Range("K3").Value = 78
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 236

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 174

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 44
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 209
$ws.Range("K8").Value = 484
$ws.Range("K11").Value = 157
$ws.Range("K19").Value = 201
$ws.Range("K20").Value = 151
$ws.Range("K23").Value = 65
$ws.Range("J29").Value = 1558
$ws.Range("K29").Value = 352
$ws.Range("K31").Value = 81
$ws.Range("K33").Value = 288
$ws.Range("J36").Value = 404
$ws.Range("K37").Value = 236
$ws.Range("K39").Value = 10
$ws.Range("K42").Value = 249
$ws.Range("K44").Value = 70
$ws.Range("K47").Value = 44
$ws.Range("K48").Value = 89
$ws.Range("K49").Value = 49
$ws.Range("K50").Value = 45
$ws.Range("K53").Value = 107
$ws.Range("K54").Value = 129
$ws.Range("F63").Value = 192
$ws.Range("J63").Value = 95
$ws.Range("K63").Value = 23
$ws.Range("K64").Value = 46
$ws.Range("K65").Value = 174
$ws.Range("K67").Value = 276
$ws.Range("K71").Value = 20
$ws.Range("K73").Value = 71
$ws.Range("K75").Value = 30
$ws.Range("K76").Value = 106
$ws.Range("K77").Value = 51
$ws.Range("K78").Value = 94
$ws.Range("K79").Value = 190
$ws.Range("K84").Value = 51
$ws.Range("K85").Value = 349
$ws.Range("K88").Value = 93
$ws.Range("K89").Value = 99
$ws.Range("K90").Value = 58
$ws.Range("K94").Value = 85
$ws.Range("K95").Value = 112
$ws.Range("K96").Value = 98
$ws.Range("K97").Value = 64
$ws.Range("K99").Value = 131
$ws.Range("F101").Value = 24100
$ws.Range("J101").Value = 29276
$ws.Range("K101").Value = 7191

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 29
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 79
$ws.Range("K3").Value = 87
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 276

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 17
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K2").Value = 4
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 44
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 120
$ws.Range("J4").Value = 85
$ws.Range("K4").Value = 18
$ws.Range("J7").Value = 1558
$ws.Range("K7").Value = 352

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 65
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K3").Value = 9
$ws.Range("K6").Value = 29

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 73
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 249

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 28
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 70
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 190

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 404

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K5").Value = 12
$ws.Range("K7").Value = 209

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 24
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 16
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("K2").Value = 4
$ws.Range("K6").Value = 10

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 48
$ws.Range("K3").Value = 39
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 12
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 20
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 118
$ws.Range("K7").Value = 349

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 51
